$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (2008年, 2009年) -- this shifts all
# subsequent rows (2010年..2020年) up by two rows.
$ws.Range("A2:A3").EntireRow.Delete() | Out-Null

# Append the new 2021年 row of data at the end (row 13 after the deletion).
# Copy the formatting from the row above (2020年) so the new year label
# keeps the same style (centered, bold border) as the rest of column A.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$newRowValues = @("2021年", 129540, 14596, 776316, 111651, 353779, 243625, 421858, 130057, 299232, 385478, 1716450, 1263462, 733098, 510756, 1972142, 716286, 64807, 1135886, 88227, 349226, 388338)

for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item(13, $i + 1).Value = $newRowValues[$i]
}
